$wb = $excel.ActiveWorkbook

# The localization file "918a4ba8-a289-4a33-ad17-42268040a38d.md" has been handed
# back for both zh-cn and de-de. Update status + handback datetime accordingly.

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: update Status columns for both languages ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: update Status + Latest Handback DateTime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G3").Value = "2016-03-01 06:30:55"

# --- de-de sheet: update Status + Latest Handback DateTime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G3").Value = "2016-03-01 06:31:19"
